$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Attribute","Y","R²","RMSE","Offset","Slope")
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$data = @(
    @("SST", "Referência", 0.8199335940901312, 1.127220969087652, 2.520106522025437, 0.8199335940901306),
    @("SST", "Predição", 0.7346559650931279, 1.368350563738811, 2.888446474693278, 0.7940417253792327),
    @("SST", "Validação", 0.7655841729998452, 1.018845322461626, 1.140579832143205, 0.9184848592463797),
    @("PH", "Referência", 0.6194159694246284, 0.190269261597771, 1.250035860105419, 0.6194159694246288),
    @("PH", "Predição", 0.4218790728609437, 0.2345053188198085, 1.638481421637806, 0.500726272970014),
    @("PH", "Validação", 0.266817439306803, 0.2256539394122687, 1.310329444091793, 0.6113653378733529),
    @("AT", "Referência", 0.6414071220225402, 0.3511647500287419, 0.4026198869844502, 0.6414071220225416),
    @("AT", "Predição", 0.4217624645620563, 0.4459264415167446, 0.5280357467405109, 0.5332948426516071),
    @("AT", "Validação", 0.686274151730107, 0.244982744917667, -0.04118963980369106, 0.9987623197514067),
    @("FIRMEZA (N)", "Referência", 0.6496827478533534, 60.82710683345287, 179.7162855527577, 0.6496827478533519),
    @("FIRMEZA (N)", "Predição", 0.384264866253081, 80.64238526564257, 245.1575041810685, 0.5208172050466638),
    @("FIRMEZA (N)", "Validação", 0.5297024350879151, 52.98044118784288, 52.66918871666712, 0.8818941314237686),
    @("UBS (%)", "Referência", 0.6474430404511626, 1.910040297746731, 5.407883862135076, 0.6474430404511622),
    @("UBS (%)", "Predição", 0.5981015966807981, 2.039323032879718, 5.813119203718665, 0.6206258724684867),
    @("UBS (%)", "Validação", 0.5953310123975417, 1.378569044227535, 3.397207338705432, 0.787222122260882)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
